$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.624.10'
$ws.Range("E2").Value = '  +0.89%  '
$ws.Range("D3").Value = '1.820.26'
$ws.Range("E3").Value = '  +1.61%  '
$ws.Range("D4").Value = '''1.009'
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''1.008'
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D6").Value = '''305.90'
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("D7").Value = '''0.4674'
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D8").Value = '''0.3600'
$ws.Range("D8").Style = $ws.Range("B8").Style
$ws.Range("D9").Value = '''46.24'
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").Value = '''0.07128'
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").Value = '''0.9008'
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = '  +2.88%  '
$ws.Range("D12").Value = '''0.07803'
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '''19.42'
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '1.801.01'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '''5.248'
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '''6.332'
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '''87.48'
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = '  +2.96%  '
$ws.Range("D18").Value = '''1.010'
$ws.Range("D18").Style = $ws.Range("B18").Style
$ws.Range("E18").Value = '  +0.11%  '
$ws.Range("D19").Value = '''0.000008568'
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = '  +0.43%  '
$ws.Range("D20").Value = '''1.007'
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").Value = '26.661.31'
$ws.Range("E21").Value = '  +0.90%  '
$ws.Range("D22").Value = '''14.19'
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = '  -0.32%  '
$ws.Range("D23").Value = '''5.010'
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").Value = '''10.57'
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = '  +0.60%  '
$ws.Range("D25").Value = '''1.931'
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("D26").Value = '''151.93'
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").Value = '''17.89'
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = '''1.984'
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = '  -2.90%  '
$ws.Range("D29").Value = '''113.66'
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = '  +1.24%  '
$ws.Range("D30").Value = '''4.801'
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("D31").Value = '''0.08791'
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = '  +1.43%  '
$ws.Range("D32").Value = '''3.147'
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = '  +2.90%  '
$ws.Range("D33").Value = '''2.768'
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("D34").Value = '''0.7315'
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("D36").Value = '''1.121'
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = '  +1.07%  '
$ws.Range("D37").Value = '''1.077'
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = '''0.01929'
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = '  -0.87%  '
$ws.Range("D39").Value = '''2.918'
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value = '  +1.77%  '
$ws.Range("D40").Value = '''0.05113'
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = '  +0.18%  '
$ws.Range("D41").Value = '''0.5063'
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = '  -3.14%  '
$ws.Range("D42").Value = '''6.813'
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = '  -1.17%  '
$ws.Range("D43").Value = '''0.1497'
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = '  -1.02%  '
$ws.Range("D44").Value = '''7.996'
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").Value = '''0.4680'
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = '  -0.67%  '
$ws.Range("D46").Value = '''1.008'
$ws.Range("D46").Style = $ws.Range("B46").Style
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("D47").Value = '''10.06'
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").Value = '''98.93'
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("D49").Value = '''1.560'
$ws.Range("D49").Style = $ws.Range("B49").Style
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("D50").Value = '''0.06016'
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = '  +0.66%  '
$ws.Range("D51").Value = '''63.62'
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = '  -0.53%  '
